$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds the "last changed" date for every logging
# notification row. The whole column (rows 2-494) is being bumped from
# 2023-09-10 (serial 45179) to 2023-09-11 (serial 45180).
$ws.Range("C2:C494").Value = 45180
